$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

# Row 4
$ws.Range("C4").Value = 3
$ws.Range("E4").Value = 5

# Row 5
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 10
$ws.Range("F5").Value = 2

# Row 6
$ws.Range("E6").Value = 0
